$wb = $excel.ActiveWorkbook

# --- mainTimeline: add the PR xref (do this before the studyDesign edit so new
#     shared strings are appended in the same order as the authored workbook).
#     The leading apostrophe keeps the cell's existing text/quote-prefix style
#     (it is not stored as part of the cell's text). ---
$wsMain = $wb.Worksheets.Item("mainTimeline")
$wsMain.Cells.Item(11, 3).Value = "'PR:profile1"
$wsMain.Range("C12").Select()

# --- studyDesign: update selection cursor and the profile reference label ---
$wsDesign = $wb.Worksheets.Item("studyDesign")
$wsDesign.Cells.Item(8, 2).Value = "profile1"
$wsDesign.Range("C23").Select()

# --- profileTimeline -> profile1: rename the sheet and move the selection cursor ---
$wsProfile = $wb.Worksheets.Item("profileTimeline")
$wsProfile.Name = "profile1"
$wsProfile.Range("C45").Select()
